$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Product Name" column (C) for
# the new "Material Type" column.
$ws.Columns("C").Insert()

# Header for the new column.
$ws.Range("C1").Value = "Material Type"

# Data rows for the new column.
$ws.Range("C4").Value = "RNA:Total RNA"
$ws.Range("C5").Value = "RNA:Total RNA"
$ws.Range("C6").Value = "RNA:Total RNA"
$ws.Range("C7").Value = "RNA:Total RNA"

# The new data cells use a plain (non-themed) Calibri font.
$ws.Range("C4:C7").Font.Name = "Calibri"

# Selection moved to C7 in the saved file.
$ws.Range("C7").Select()
